# "Generate Report for Handback"
#
# The localization round-trip finished: both files ("34993c0e...md" and
# "fd49591c...md") have been handed back and are in sync with en-US. This
# script records that on the per-language report sheets (zh-cn, de-de):
#   - Status (Overview + per-language "Status" column) flips from
#     "In Translation" to "Handed back: in sync with en-US".
#   - "Latest Target File" (col I) now links to the source .md file.
#   - "Latest Handback File" (col J) is filled in with the generated
#     per-language .xlf handback file name.
#   - "Latest Handback DateTime" (col K) is stamped with the handback time.
#   - A few report columns are widened so the new long file names are
#     readable.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f0a2de0ebf4fed3768bef992fe19781b86f531e8/e2e"

$newStatus = "Handed back: in sync with en-US"

$file1Md  = "34993c0e-405d-4214-9d97-07a97e3e8ca2.md"
$file2Md  = "fd49591c-0f33-42cb-9fc5-4d18d27504ba.md"

$file1ZhCnXlf = "34993c0e-405d-4214-9d97-07a97e3e8ca2.c0a0ed21ec3225906d6aa5cc7f151138fb9362a4.zh-cn.xlf"
$file2ZhCnXlf = "fd49591c-0f33-42cb-9fc5-4d18d27504ba.22ded12f3844b568e32595dbfa60a7d574af4227.zh-cn.xlf"
$file1DeDeXlf = "34993c0e-405d-4214-9d97-07a97e3e8ca2.c0a0ed21ec3225906d6aa5cc7f151138fb9362a4.de-de.xlf"
$file2DeDeXlf = "fd49591c-0f33-42cb-9fc5-4d18d27504ba.22ded12f3844b568e32595dbfa60a7d574af4227.de-de.xlf"

$zhCnHandbackTime = "2016-09-03 04:28:13"
$deDeHandbackTime = "2016-09-03 04:28:19"

# ---------------------------------------------------------------------
# Overview sheet: status cells just pick up the shared-string text
# change automatically; widen the two status columns (zh-cn / de-de).
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(2, 5).Value = $newStatus
$wsOverview.Cells.Item(2, 6).Value = $newStatus
$wsOverview.Cells.Item(3, 5).Value = $newStatus
$wsOverview.Cells.Item(3, 6).Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(2, 3).Value = $newStatus
$wsZhCn.Cells.Item(3, 3).Value = $newStatus

# Row 2 (34993c0e...)
$wsZhCn.Cells.Item(2, 9).Value = $file1Md
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(2, 9), "$baseUrl/$file1Md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $file1Md)
$wsZhCn.Cells.Item(2, 10).Value = $file1ZhCnXlf
$wsZhCn.Cells.Item(2, 11).Value = $zhCnHandbackTime

# Row 3 (fd49591c...)
$wsZhCn.Cells.Item(3, 9).Value = $file2Md
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(3, 9), "$baseUrl/$file2Md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $file2Md)
$wsZhCn.Cells.Item(3, 10).Value = $file2ZhCnXlf
$wsZhCn.Cells.Item(3, 11).Value = $zhCnHandbackTime

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(2, 3).Value = $newStatus
$wsDeDe.Cells.Item(3, 3).Value = $newStatus

# Row 2 (34993c0e...)
$wsDeDe.Cells.Item(2, 9).Value = $file1Md
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(2, 9), "$baseUrl/$file1Md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $file1Md)
$wsDeDe.Cells.Item(2, 10).Value = $file1DeDeXlf
$wsDeDe.Cells.Item(2, 11).Value = $deDeHandbackTime

# Row 3 (fd49591c...)
$wsDeDe.Cells.Item(3, 9).Value = $file2Md
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(3, 9), "$baseUrl/$file2Md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $file2Md)
$wsDeDe.Cells.Item(3, 10).Value = $file2DeDeXlf
$wsDeDe.Cells.Item(3, 11).Value = $deDeHandbackTime

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
